$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '35.259.25'
$ws.Range('E2').Value = '  +2.04%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.842.51'
$ws.Range('E3').Value = '  +1.90%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '228.09'
$ws.Range('E5').Value = '  +1.06%  '

$ws.Range('E6').Value = '  +2.55%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '43.32'
$ws.Range('E8').Value = '  +15.59%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.305'
$ws.Range('E9').Value = '  +4.72%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0693'
$ws.Range('E10').Value = '  +1.81%  '

$ws.Range('E11').Value = '  +3.43%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.109.02'
$ws.Range('E12').Value = '  +1.86%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '11.58'
$ws.Range('E13').Value = '  +2.40%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.839.18'
$ws.Range('E14').Value = '  +1.64%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.73'
$ws.Range('E15').Value = '  +6.87%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.660'
$ws.Range('E16').Value = '  +4.38%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '35.210.97'
$ws.Range('E17').Value = '  +2.06%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '70.01'
$ws.Range('E18').Value = '  +1.97%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '245.79'
$ws.Range('E19').Value = '  +0.91%  '

$ws.Range('E20').Value = '  +2.38%  '

$ws.Range('E21').Value = '  +8.20%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.71'
$ws.Range('E22').Value = '  +13.92%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.04%  '

$ws.Range('E24').Value = '  -1.08%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '172.38'
$ws.Range('E25').Value = '  +0.19%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.92'
$ws.Range('E26').Value = '  +1.12%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '17.90'
$ws.Range('E27').Value = '  +3.47%  '

$ws.Range('E28').Value = '  +1.82%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '3.646.37'
$ws.Range('E29').Value = '  +50.08%  '

$ws.Range('E30').Value = '  -0.08%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.33'
$ws.Range('E31').Value = '  +7.86%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.94'
$ws.Range('E32').Value = '  +3.34%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.05'
$ws.Range('E33').Value = '  +3.39%  '

$ws.Range('E34').Value = '  +3.80%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.88'
$ws.Range('E35').Value = '  +3.89%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.673'
$ws.Range('E36').Value = '  +3.00%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '90.39'
$ws.Range('E37').Value = '  +11.95%  '

$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.08'
$ws.Range('E38').Value = '  +1.35%  '

$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.342.18'
$ws.Range('E39').Value = '  -1.59%  '

$ws.Range('E40').Value = '  +8.92%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.43'
$ws.Range('E41').Value = '  +2.84%  '

$ws.Range('E42').Value = '  +3.59%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '14.88'
$ws.Range('E43').Value = '  +8.67%  '

$ws.Range('E44').Value = '  +6.54%  '

$ws.Range('E45').Value = '  +1.02%  '

$ws.Range('E46').Value = '  +1.74%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0520'
$ws.Range('E47').Value = '  +3.65%  '

$ws.Range('E48').Value = '  +4.48%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.009.78'
$ws.Range('E49').Value = '  +1.98%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '104.75'
$ws.Range('E50').Value = '  +2.15%  '
